# Working on NOT NULL constraint
# - Rename the *_ID primary-key columns to CamelCase "...ID" (no underscore)
#   across Schema, Member, Role and Tour sheets.
# - Add a new "MemberRole" junction table to the Schema sheet (between Tour
#   and TourAge), and rename TourAge's foreign-key columns to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Schema sheet
# ---------------------------------------------------------------------
$schema = $wb.Worksheets.Item("Schema")

# Rename the primary-key id columns for the existing base tables.
$schema.Cells.Item(2, 2).Value = "MemberID"
$schema.Cells.Item(4, 2).Value = "RoleID"
$schema.Cells.Item(6, 2).Value = "TourID"

# Insert 3 new rows for the "MemberRole" table right after the "Tour" rows
# (before the old "TourAge" block, which starts at row 8).
$null = $schema.Rows.Item(8).Resize(3).Insert()

$schema.Cells.Item(8, 1).Value = "MemberRole"
$schema.Cells.Item(8, 2).Value = "MemberRoleID"
$schema.Cells.Item(8, 3).Value = "int"
$schema.Cells.Item(8, 4).Value = "NA"
$schema.Cells.Item(8, 5).Value = "NA"

$schema.Cells.Item(9, 1).Value = "MemberRole"
$schema.Cells.Item(9, 2).Value = "Member"
$schema.Cells.Item(9, 3).Value = "int"
$schema.Cells.Item(9, 4).Value = "Member"
$schema.Cells.Item(9, 5).Value = "MemberID"

$schema.Cells.Item(10, 1).Value = "MemberRole"
$schema.Cells.Item(10, 2).Value = "Role"
$schema.Cells.Item(10, 3).Value = "int"
$schema.Cells.Item(10, 4).Value = "Role"
$schema.Cells.Item(10, 5).Value = "RoleID"

# The old "TourAge" rows have shifted down to 11-14; rename its id + fk cols.
$schema.Cells.Item(11, 2).Value = "TourAgeID"
$schema.Cells.Item(12, 2).Value = "MemberID"
$schema.Cells.Item(12, 5).Value = "MemberID"
$schema.Cells.Item(13, 2).Value = "TourID"
$schema.Cells.Item(13, 5).Value = "TourID"

# Column B now holds longer names ("MemberRoleID") - widen it to fit.
$schema.Columns.Item(2).ColumnWidth = 14.1640625

$null = $schema.Range("E13").Select()

# ---------------------------------------------------------------------
# Member sheet
# ---------------------------------------------------------------------
$member = $wb.Worksheets.Item("Member")
$member.Cells.Item(1, 1).Value = "MemberID"

# ---------------------------------------------------------------------
# Role sheet
# ---------------------------------------------------------------------
$role = $wb.Worksheets.Item("Role")
$role.Cells.Item(1, 1).Value = "RoleID"

# ---------------------------------------------------------------------
# Tour sheet
# ---------------------------------------------------------------------
$tour = $wb.Worksheets.Item("Tour")
$tour.Cells.Item(1, 1).Value = "TourID"
$null = $tour.Range("J23").Select()

Write-Output "Schema updated: MemberRole table added, id columns renamed."
